$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two input cells; dependent formula cells (D6, D7, D8, D9, E8, E9)
# recalculate automatically.
$ws.Range("D3").Value = 422421.16
$ws.Range("D5").Value = 147677.48000000001

# Move the active selection to D5 (as last edited/selected cell).
$ws.Range("D5").Select()
